{"js": "// Update the division-problem answers in the worksheet table.\n// Each cell holds a unique \"a\u00f7b=c, d\" string, so we can safely locate\n// and replace each one independently via Body.search().\nconst replacements = [\n  [\"20\u00f78=2, 4\", \"46\u00f77=6, 4\"],\n  [\"77\u00f77=11, 0\", \"26\u00f75=5, 1\"],\n  [\"77\u00f74=19, 1\", \"11\u00f72=5, 1\"],\n  [\"95\u00f76=15, 5\", \"87\u00f74=21, 3\"],\n  [\"54\u00f72=27, 0\", \"95\u00f77=13, 4\"],\n  [\"52\u00f78=6, 4\", \"38\u00f74=9, 2\"],\n  [\"16\u00f75=3, 1\", \"86\u00f74=21, 2\"],\n  [\"80\u00f76=13, 2\", \"88\u00f73=29, 1\"],\n  [\"72\u00f76=12, 0\", \"12\u00f74=3, 0\"],\n  [\"57\u00f74=14, 1\", \"90\u00f74=22, 2\"],\n  [\"72\u00f78=9, 0\", \"39\u00f79=4, 3\"],\n  [\"32\u00f78=4, 0\", \"15\u00f75=3, 0\"],\n  [\"88\u00f78=11, 0\", \"23\u00f79=2, 5\"],\n  [\"37\u00f74=9, 1\", \"17\u00f77=2, 3\"],\n  [\"47\u00f78=5, 7\", \"35\u00f79=3, 8\"],\n  [\"18\u00f79=2, 0\", \"59\u00f74=14, 3\"],\n  [\"14\u00f76=2, 2\", \"56\u00f76=9, 2\"],\n  [\"45\u00f79=5, 0\", \"83\u00f72=41, 1\"],\n  [\"22\u00f78=2, 6\", \"65\u00f72=32, 1\"],\n  [\"72\u00f77=10, 2\", \"59\u00f76=9, 5\"],\n  [\"81\u00f79=9, 0\", \"14\u00f78=1, 6\"],\n  [\"99\u00f77=14, 1\", \"90\u00f75=18, 0\"],\n  [\"94\u00f76=15, 4\", \"67\u00f73=22, 1\"],\n  [\"64\u00f74=16, 0\", \"27\u00f75=5, 2\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the division-problem answers in the worksheet table.\n# Each cell holds a unique \"a\u00f7b=c, d\" string, so Find/Replace on the\n# whole document is safe and unambiguous for every pair below.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"20\u00f78=2, 4\", \"46\u00f77=6, 4\"),\n    @(\"77\u00f77=11, 0\", \"26\u00f75=5, 1\"),\n    @(\"77\u00f74=19, 1\", \"11\u00f72=5, 1\"),\n    @(\"95\u00f76=15, 5\", \"87\u00f74=21, 3\"),\n    @(\"54\u00f72=27, 0\", \"95\u00f77=13, 4\"),\n    @(\"52\u00f78=6, 4\", \"38\u00f74=9, 2\"),\n    @(\"16\u00f75=3, 1\", \"86\u00f74=21, 2\"),\n    @(\"80\u00f76=13, 2\", \"88\u00f73=29, 1\"),\n    @(\"72\u00f76=12, 0\", \"12\u00f74=3, 0\"),\n    @(\"57\u00f74=14, 1\", \"90\u00f74=22, 2\"),\n    @(\"72\u00f78=9, 0\", \"39\u00f79=4, 3\"),\n    @(\"32\u00f78=4, 0\", \"15\u00f75=3, 0\"),\n    @(\"88\u00f78=11, 0\", \"23\u00f79=2, 5\"),\n    @(\"37\u00f74=9, 1\", \"17\u00f77=2, 3\"),\n    @(\"47\u00f78=5, 7\", \"35\u00f79=3, 8\"),\n    @(\"18\u00f79=2, 0\", \"59\u00f74=14, 3\"),\n    @(\"14\u00f76=2, 2\", \"56\u00f76=9, 2\"),\n    @(\"45\u00f79=5, 0\", \"83\u00f72=41, 1\"),\n    @(\"22\u00f78=2, 6\", \"65\u00f72=32, 1\"),\n    @(\"72\u00f77=10, 2\", \"59\u00f76=9, 5\"),\n    @(\"81\u00f79=9, 0\", \"14\u00f78=1, 6\"),\n    @(\"99\u00f77=14, 1\", \"90\u00f75=18, 0\"),\n    @(\"94\u00f76=15, 4\", \"67\u00f73=22, 1\"),\n    @(\"64\u00f74=16, 0\", \"27\u00f75=5, 2\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $rng = $d.Content\n    $found = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        Write-Output \"NOT FOUND: $oldText\"\n    }\n}\n"}
